$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("C2").Value = "Kategoria"
$ws.Range("D2").Value = "Opis"

# Move the selection to match the final state (side effect of editing)
$ws.Range("G7").Select()
